$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:D20").Clear()

    $ws.Range("A1").Value = "Target"
    $ws.Range("B1").Value = "Start"
    $ws.Range("C1").Value = "End"
    $ws.Range("D1").Value = "Torus"
    $ws.Range("A2").Value = "LD2"
    $ws.Range("B2").Value = 18305
    $ws.Range("C2").Value = 18312
    $ws.Range("D2").Value = "IB"
    $ws.Range("A3").Value = "Empty "
    $ws.Range("B3").Value = 18316
    $ws.Range("C3").Value = 18316
    $ws.Range("D3").Value = "Zero-field "
    $ws.Range("A4").Value = "LD2"
    $ws.Range("B4").Value = 18318
    $ws.Range("C4").Value = 18336
    $ws.Range("D4").Value = "IB"
    $ws.Range("A5").Value = "CxC"
    $ws.Range("B5").Value = 18339
    $ws.Range("C5").Value = 18346
    $ws.Range("D5").Value = "IB"
    $ws.Range("A6").Value = "CuSn"
    $ws.Range("B6").Value = 18347
    $ws.Range("C6").Value = 18368
    $ws.Range("D6").Value = "IB"
    $ws.Range("A7").Value = "CxC"
    $ws.Range("B7").Value = 18369
    $ws.Range("C7").Value = 18371
    $ws.Range("D7").Value = "IB"
    $ws.Range("A8").Value = "CuSn"
    $ws.Range("B8").Value = 18372
    $ws.Range("C8").Value = 18394
    $ws.Range("D8").Value = "IB"
    $ws.Range("A9").Value = "Empty "
    $ws.Range("B9").Value = 18399
    $ws.Range("C9").Value = 18399
    $ws.Range("D9").Value = "IB"
    $ws.Range("A10").Value = "CxC"
    $ws.Range("B10").Value = 18400
    $ws.Range("C10").Value = 18401
    $ws.Range("D10").Value = "IB"
    $ws.Range("A11").Value = "LD2"
    $ws.Range("B11").Value = 18419
    $ws.Range("C11").Value = 18439
    $ws.Range("D11").Value = "OB"
    $ws.Range("A12").Value = "CxC"
    $ws.Range("B12").Value = 18440
    $ws.Range("C12").Value = 18524
    $ws.Range("D12").Value = "OB"
    $ws.Range("A13").Value = "LD2"
    $ws.Range("B13").Value = 18528
    $ws.Range("C13").Value = 18559
    $ws.Range("D13").Value = "OB"
    $ws.Range("A14").Value = "CuSn"
    $ws.Range("B14").Value = 18560
    $ws.Range("C14").Value = 18642
    $ws.Range("D14").Value = "OB"
    $ws.Range("A15").Value = "LD2"
    $ws.Range("B15").Value = 18644
    $ws.Range("C15").Value = 18656
    $ws.Range("D15").Value = "OB"
    $ws.Range("A16").Value = "CuSn"
    $ws.Range("B16").Value = 18660
    $ws.Range("C16").Value = 18755
    $ws.Range("D16").Value = "OB"
    $ws.Range("A17").Value = "CxC"
    $ws.Range("B17").Value = 18756
    $ws.Range("C17").Value = 18762
    $ws.Range("D17").Value = "OB"
    $ws.Range("A18").Value = "LD2"
    $ws.Range("B18").Value = 18764
    $ws.Range("C18").Value = 18790
    $ws.Range("D18").Value = "OB"
    $ws.Range("A19").Value = "CxC"
    $ws.Range("B19").Value = 18796
    $ws.Range("C19").Value = 18850
    $ws.Range("D19").Value = "OB"
    $ws.Range("A20").Value = "LD2"
    $ws.Range("B20").Value = 18851
    $ws.Range("C20").Value = 18873
    $ws.Range("D20").Value = "OB"
    $ws.Range("A21").Value = "CuSn"
    $ws.Range("B21").Value = 18874
    $ws.Range("C21").Value = 18966
    $ws.Range("D21").Value = "OB"
    $ws.Range("A22").Value = "LD2"
    $ws.Range("B22").Value = 19021
    $ws.Range("C22").Value = 19058
    $ws.Range("D22").Value = "OB"
    $ws.Range("A23").Value = "Empty "
    $ws.Range("B23").Value = 19060
    $ws.Range("C23").Value = 19060
    $ws.Range("D23").Value = "OB"
    $ws.Range("A24").Value = "CuSn"
    $ws.Range("B24").Value = 19061
    $ws.Range("C24").Value = 19131
    $ws.Range("D24").Value = "OB"

$rng = $ws.Range("A1:D24")
$rng.Font.Size = 11
$rng.Font.Name = "Arial"

$ws.Range("E5").Select()
